# reviewdb.xlsx - "Add files via upload"
#
# Net effect of the commit (once the shared-string shuffle noise in the
# diff is resolved against the cell indices that reference it):
#   - Sheet1!G4 ("blue"/recovery-confirmation column) changes from
#     "confirm" to "no".
#   - The rows for G5/G6 keep referring to the same text ("confirm"); the
#     shared-strings table just got reordered by Excel when it rewrote the
#     file, so there is no real content change there.
#   - The saved selection/view state moves from C11:D12 (active cell D11)
#     to a single-cell selection on G5, with the viewport scrolled one
#     column to the right (topLeftCell B1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Only real content change in the workbook.
$ws.Range("G4").Value = "no"

# Restore/ensure the other "blue" cells touched by the diff keep their
# original text (no-ops content-wise, just guards against drift).
$ws.Range("G5").Value = "confirm"
$ws.Range("G6").Value = "confirm"

# Match the new saved view state: active cell / selection on G5.
$ws.Range("G5").Select()
